$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 452.25
$ws.Range("I58").Value = 452.25
$ws.Range("K58").Value = 1356.75
$ws.Range("M58").Value = -1206.75
$ws.Range("H62").Value = 16172.579
$ws.Range("I62").Value = 15477.9
$ws.Range("K62").Value = 15477.9
$ws.Range("M62").Value = -14853.9
$ws.Range("H65").Value = 16172.579
$ws.Range("I65").Value = 15477.9
$ws.Range("K65").Value = 77389.5
$ws.Range("M65").Value = -74269.5
$ws.Range("H100").Value = 6471.706
$ws.Range("I100").Value = 1444.3334
$ws.Range("J100").Value = 12127.5
$ws.Range("K100").Value = 1444.3334
$ws.Range("L100").Value = 12127.5
$ws.Range("M100").Value = -903.3334
$ws.Range("N100").Value = -13209.5
$ws.Range("H136").Value = 88999.5
$ws.Range("J136").Value = 88999.5
$ws.Range("L136").Value = 88999.5
$ws.Range("N136").Value = -99199.5
$ws.Range("H137").Value = 422677.2
$ws.Range("I137").Value = 662678.5
$ws.Range("J137").Value = 2674.8333
$ws.Range("K137").Value = 1988035.5
$ws.Range("L137").Value = 8024.499899999999
$ws.Range("M137").Value = -1985485.5
$ws.Range("N137").Value = -13124.4999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3440677.5
$ws.Range("I32").Value = 4117984.8
$ws.Range("K32").Value = 4117984.8
$ws.Range("M32").Value = -4117697.8
$ws.Range("H36").Value = 7725.75
$ws.Range("I36").Value = 9001.5
$ws.Range("K36").Value = 9001.5
$ws.Range("M36").Value = -8655.5
$ws.Range("H121").Value = 17998
$ws.Range("J121").Value = 17998
$ws.Range("L121").Value = 17998
$ws.Range("N121").Value = -21492

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3499
$ws.Range("I86").Value = 3998
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3998
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -2875
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3499
$ws.Range("I89").Value = 3998
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 19990
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -14374
$ws.Range("N89").Value = -26232
$ws.Range("H94").Value = 386.5
$ws.Range("I94").Value = 386.5
$ws.Range("K94").Value = 386.5
$ws.Range("M94").Value = 64.5
$ws.Range("H135").Value = 83584
$ws.Range("J135").Value = 83584
$ws.Range("L135").Value = 83584
$ws.Range("N135").Value = -93724

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5378.0205
$ws.Range("I31").Value = 2045.3846
$ws.Range("J31").Value = 9145.348
$ws.Range("K31").Value = 2045.3846
$ws.Range("L31").Value = 9145.348
$ws.Range("M31").Value = -1750.3846
$ws.Range("N31").Value = -9735.348
$ws.Range("H34").Value = 5378.0205
$ws.Range("I34").Value = 2045.3846
$ws.Range("J34").Value = 9145.348
$ws.Range("K34").Value = 2045.3846
$ws.Range("L34").Value = 9145.348
$ws.Range("M34").Value = -1843.3846
$ws.Range("N34").Value = -9549.348
$ws.Range("H56").Value = 1500
$ws.Range("I56").Value = 1500
$ws.Range("K56").Value = 1500
$ws.Range("M56").Value = -655
$ws.Range("H88").Value = 14820.375
$ws.Range("J88").Value = 14820.375
$ws.Range("L88").Value = 14820.375
$ws.Range("N88").Value = -15632.375
$ws.Range("H91").Value = 14820.375
$ws.Range("J91").Value = 14820.375
$ws.Range("L91").Value = 14820.375
$ws.Range("N91").Value = -17628.375
$ws.Range("H99").Value = 5667.8667
$ws.Range("I99").Value = 5316.769
$ws.Range("K99").Value = 5316.769
$ws.Range("M99").Value = -3818.769
$ws.Range("H105").Value = 29461.916
$ws.Range("I105").Value = 32095.092
$ws.Range("J105").Value = 497
$ws.Range("K105").Value = 32095.092
$ws.Range("L105").Value = 497
$ws.Range("M105").Value = -30348.092
$ws.Range("N105").Value = -3991
$ws.Range("H126").Value = 5667.8667
$ws.Range("I126").Value = 5316.769
$ws.Range("K126").Value = 15950.307
$ws.Range("M126").Value = -13480.307
$ws.Range("H140").Value = 92166.59
$ws.Range("J140").Value = 92166.59
$ws.Range("L140").Value = 92166.59
$ws.Range("N140").Value = -102526.59

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 602.7143
$ws.Range("I9").Value = 602.7143
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1808.1429
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -1584.1429
$ws.Range("N9").ClearContents()
$ws.Range("H34").Value = 1825.9286
$ws.Range("J34").Value = 5047.25
$ws.Range("L34").Value = 15141.75
$ws.Range("N34").Value = -15309.75
$ws.Range("H39").Value = 2784.3333
$ws.Range("J39").Value = 2940.6
$ws.Range("L39").Value = 8821.799999999999
$ws.Range("N39").Value = -9409.799999999999
$ws.Range("H55").Value = 2054.2
$ws.Range("J55").Value = 2133.3333
$ws.Range("L55").Value = 6399.999899999999
$ws.Range("N55").Value = -6753.999899999999
$ws.Range("H70").Value = 4382.5
$ws.Range("I70").Value = 1930.2858
$ws.Range("K70").Value = 5790.857400000001
$ws.Range("M70").Value = -5475.857400000001
$ws.Range("H73").Value = 4382.5
$ws.Range("I73").Value = 1930.2858
$ws.Range("K73").Value = 5790.857400000001
$ws.Range("M73").Value = -4698.857400000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6761.6924
$ws.Range("I80").Value = 5375.5
$ws.Range("J80").Value = 8979.6
$ws.Range("K80").Value = 5375.5
$ws.Range("L80").Value = 8979.6
$ws.Range("M80").Value = -4377.5
$ws.Range("N80").Value = -10975.6
$ws.Range("H83").Value = 6761.6924
$ws.Range("I83").Value = 5375.5
$ws.Range("J83").Value = 8979.6
$ws.Range("K83").Value = 26877.5
$ws.Range("L83").Value = 44898
$ws.Range("M83").Value = -21885.5
$ws.Range("N83").Value = -54882
$ws.Range("H102").Value = 3559.2632
$ws.Range("I102").Value = 1663.4
$ws.Range("K102").Value = 1663.4
$ws.Range("M102").Value = -41.40000000000009

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H61").Value = 5623.5
$ws.Range("I61").Value = 3997.25
$ws.Range("K61").Value = 3997.25
$ws.Range("M61").Value = -3795.25
$ws.Range("H93").Value = 2921.5454
$ws.Range("I93").Value = 3297.5
$ws.Range("K93").Value = 3297.5
$ws.Range("M93").Value = -2049.5
$ws.Range("H113").Value = 5623.5
$ws.Range("I113").Value = 3997.25
$ws.Range("K113").Value = 3997.25
$ws.Range("M113").Value = -1827.25
$ws.Range("H122").Value = 4623.3184
$ws.Range("J122").Value = 5022.5713
$ws.Range("L122").Value = 15067.7139
$ws.Range("N122").Value = -19967.7139

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 17500
$ws.Range("I49").Value = 15000
$ws.Range("K49").Value = 15000
$ws.Range("M49").Value = -14770
$ws.Range("H113").Value = 2831.2222
$ws.Range("I113").Value = 740.0833
$ws.Range("J113").Value = 7013.5
$ws.Range("K113").Value = 2220.2499
$ws.Range("L113").Value = 21040.5
$ws.Range("M113").Value = -50.2498999999998
$ws.Range("N113").Value = -25380.5
$ws.Range("H132").Value = 4277528
$ws.Range("I132").Value = 7940634.5
$ws.Range("J132").Value = 3903.8333
$ws.Range("K132").Value = 23821903.5
$ws.Range("L132").Value = 11711.4999
$ws.Range("M132").Value = -23819373.5
$ws.Range("N132").Value = -16771.4999
